$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 71, pushing the existing rows 71-94 down to 73-96.
$ws.Rows.Item(71).Insert()
$ws.Rows.Item(71).Insert()

# --- New row 71: weekly record for "Primera" quality ---
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44463
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112021
$ws.Cells.Item(71, 7).Value = "Ají"
$ws.Cells.Item(71, 8).Value = "Inferno"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 500
$ws.Cells.Item(71, 11).Value = 41000
$ws.Cells.Item(71, 12).Value = 42000
$ws.Cells.Item(71, 13).Value = 41500
$ws.Cells.Item(71, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 3458
$ws.Cells.Item(71, 17).Value = 12
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# --- New row 72: weekly record for "Segunda" quality ---
$ws.Cells.Item(72, 1).Value = 8
$ws.Cells.Item(72, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(72, 3).Value = "Coquimbo"
$ws.Cells.Item(72, 4).Value = 44463
$ws.Cells.Item(72, 5).Value = 4
$ws.Cells.Item(72, 6).Value = 100112021
$ws.Cells.Item(72, 7).Value = "Ají"
$ws.Cells.Item(72, 8).Value = "Inferno"
$ws.Cells.Item(72, 9).Value = "Segunda"
$ws.Cells.Item(72, 10).Value = 500
$ws.Cells.Item(72, 11).Value = 34000
$ws.Cells.Item(72, 12).Value = 35000
$ws.Cells.Item(72, 13).Value = 34500
$ws.Cells.Item(72, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 2875
$ws.Cells.Item(72, 17).Value = 12
$ws.Cells.Item(72, 18).Value = "Hortaliza"
